# Auto-generated edit script: refresh market-price derived columns (H-N)
# across ALC/ARM/BSM/CRP/CUL/GSM/LTW/WVR sheets, per scheduled-runner update.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H32").Value = 1065.8
$ws.Range("J32").Value = 1199.2
$ws.Range("L32").Value = 1199.2
$ws.Range("N32").Value = -1851.2
$ws.Range("H88").Value = 2962.5
$ws.Range("I88").Value = 2962.5
$ws.Range("J88").Value = 0
$ws.Range("K88").Value = 2962.5
$ws.Range("L88").Value = 0
$ws.Range("M88").Value = -2556.5
$ws.Range("N88").ClearContents()
$ws.Range("H91").Value = 2962.5
$ws.Range("I91").Value = 2962.5
$ws.Range("J91").Value = 0
$ws.Range("K91").Value = 2962.5
$ws.Range("L91").Value = 0
$ws.Range("M91").Value = -1558.5
$ws.Range("N91").ClearContents()
$ws.Range("H107").Value = 575.1
$ws.Range("I107").Value = 478.2
$ws.Range("K107").Value = 478.2
$ws.Range("M107").Value = 1441.8
$ws.Range("H111").Value = 828.2857
$ws.Range("I111").Value = 784.25
$ws.Range("K111").Value = 2352.75
$ws.Range("M111").Value = 714.25
$ws.Range("H115").Value = 121.333336
$ws.Range("I115").Value = 121.333336
$ws.Range("K115").Value = 364.000008
$ws.Range("M115").Value = 1202.999992
$ws.Range("H116").Value = 4682.857
$ws.Range("I116").Value = 3699
$ws.Range("K116").Value = 3699
$ws.Range("M116").Value = -257
$ws.Range("H125").Value = 116315.22
$ws.Range("I125").Value = 3841.8
$ws.Range("J125").Value = 256907
$ws.Range("K125").Value = 34576.2
$ws.Range("L125").Value = 2312163
$ws.Range("M125").Value = -32116.2
$ws.Range("N125").Value = -2317083

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 599.1
$ws.Range("I2").Value = 491.22223
$ws.Range("K2").Value = 491.22223
$ws.Range("M2").Value = -378.22223
$ws.Range("H32").Value = 2163.963
$ws.Range("I32").Value = 1862.5769
$ws.Range("K32").Value = 1862.5769
$ws.Range("M32").Value = -1575.5769
$ws.Range("H74").Value = 1123.1428
$ws.Range("I74").Value = 983.3333
$ws.Range("J74").Value = 1228
$ws.Range("K74").Value = 983.3333
$ws.Range("L74").Value = 1228
$ws.Range("M74").Value = -109.3333
$ws.Range("N74").Value = -2976
$ws.Range("H77").Value = 1123.1428
$ws.Range("I77").Value = 983.3333
$ws.Range("J77").Value = 1228
$ws.Range("K77").Value = 4916.6665
$ws.Range("L77").Value = 6140
$ws.Range("M77").Value = -548.6665000000003
$ws.Range("N77").Value = -14876
$ws.Range("H110").Value = 2715.4
$ws.Range("I110").Value = 1193.5
$ws.Range("J110").Value = 4998.25
$ws.Range("K110").Value = 1193.5
$ws.Range("L110").Value = 4998.25
$ws.Range("M110").Value = 851.5
$ws.Range("N110").Value = -9088.25
$ws.Range("H116").Value = 599.1
$ws.Range("I116").Value = 491.22223
$ws.Range("K116").Value = 491.22223
$ws.Range("M116").Value = 1802.77777
$ws.Range("H122").Value = 2383.1667
$ws.Range("I122").Value = 1350
$ws.Range("J122").Value = 4449.5
$ws.Range("K122").Value = 4050
$ws.Range("L122").Value = 13348.5
$ws.Range("M122").Value = -1600
$ws.Range("N122").Value = -18248.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 599.1
$ws.Range("I3").Value = 491.22223
$ws.Range("K3").Value = 491.22223
$ws.Range("M3").Value = -377.22223
$ws.Range("H22").Value = 150
$ws.Range("I22").Value = 150
$ws.Range("J22").Value = 0
$ws.Range("K22").Value = 150
$ws.Range("L22").Value = 0
$ws.Range("M22").Value = 23
$ws.Range("N22").ClearContents()
$ws.Range("H107").Value = 1271.2142
$ws.Range("I107").Value = 1208.8182
$ws.Range("J107").Value = 1500
$ws.Range("K107").Value = 1208.8182
$ws.Range("L107").Value = 1500
$ws.Range("M107").Value = 711.1818000000001
$ws.Range("N107").Value = -5340

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 96.666664
$ws.Range("I7").Value = 95
$ws.Range("K7").Value = 95
$ws.Range("M7").Value = 18
$ws.Range("H16").Value = 996.5
$ws.Range("J16").Value = 993
$ws.Range("L16").Value = 993
$ws.Range("N16").Value = -1567
$ws.Range("H31").Value = 1731.5385
$ws.Range("I31").Value = 1731.5385
$ws.Range("K31").Value = 1731.5385
$ws.Range("M31").Value = -1436.5385
$ws.Range("H33").Value = 2697.5
$ws.Range("I33").Value = 263.33334
$ws.Range("J33").Value = 10000
$ws.Range("K33").Value = 263.33334
$ws.Range("L33").Value = 10000
$ws.Range("M33").Value = 115.66666
$ws.Range("N33").Value = -10758
$ws.Range("H34").Value = 1731.5385
$ws.Range("I34").Value = 1731.5385
$ws.Range("K34").Value = 1731.5385
$ws.Range("M34").Value = -1529.5385
$ws.Range("H107").Value = 639.1
$ws.Range("I107").Value = 602.3570999999999
$ws.Range("J107").Value = 724.8333
$ws.Range("K107").Value = 602.3570999999999
$ws.Range("L107").Value = 724.8333
$ws.Range("M107").Value = 1317.6429
$ws.Range("N107").Value = -4564.8333
$ws.Range("H113").Value = 996.5
$ws.Range("J113").Value = 993
$ws.Range("L113").Value = 993
$ws.Range("N113").Value = -5333
$ws.Range("H132").Value = 3511
$ws.Range("I132").Value = 2255.75
$ws.Range("J132").Value = 4228.2856
$ws.Range("K132").Value = 6767.25
$ws.Range("L132").Value = 12684.8568
$ws.Range("M132").Value = -4237.25
$ws.Range("N132").Value = -17744.8568

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H38").Value = 191.81818
$ws.Range("I38").Value = 175.66667
$ws.Range("J38").Value = 211.2
$ws.Range("K38").Value = 527.00001
$ws.Range("L38").Value = 633.5999999999999
$ws.Range("M38").Value = -180.00001
$ws.Range("N38").Value = -1327.6
$ws.Range("H129").Value = 550.25
$ws.Range("I129").Value = 550.25
$ws.Range("K129").Value = 1650.75
$ws.Range("M129").Value = 3349.25

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 3329.2
$ws.Range("I80").Value = 3049
$ws.Range("J80").Value = 3749.5
$ws.Range("K80").Value = 3049
$ws.Range("L80").Value = 3749.5
$ws.Range("M80").Value = -2051
$ws.Range("N80").Value = -5745.5
$ws.Range("H83").Value = 3329.2
$ws.Range("I83").Value = 3049
$ws.Range("J83").Value = 3749.5
$ws.Range("K83").Value = 15245
$ws.Range("L83").Value = 18747.5
$ws.Range("M83").Value = -10253
$ws.Range("N83").Value = -28731.5
$ws.Range("H102").Value = 1033.25
$ws.Range("I102").Value = 1033.25
$ws.Range("K102").Value = 1033.25
$ws.Range("M102").Value = 588.75

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 1045.75
$ws.Range("I22").Value = 595.5
$ws.Range("J22").Value = 1496
$ws.Range("K22").Value = 595.5
$ws.Range("L22").Value = 1496
$ws.Range("M22").Value = -300.5
$ws.Range("N22").Value = -2086
$ws.Range("H27").Value = 1045.75
$ws.Range("I27").Value = 595.5
$ws.Range("J27").Value = 1496
$ws.Range("K27").Value = 595.5
$ws.Range("L27").Value = 1496
$ws.Range("M27").Value = -488.5
$ws.Range("N27").Value = -1710
$ws.Range("H32").Value = 9707.571
$ws.Range("I32").Value = 9707.571
$ws.Range("K32").Value = 9707.571
$ws.Range("M32").Value = -9390.571
$ws.Range("H46").Value = 2000
$ws.Range("I46").Value = 2000
$ws.Range("K46").Value = 2000
$ws.Range("M46").Value = -1812
$ws.Range("H132").Value = 6174.3
$ws.Range("I132").Value = 6495.75
$ws.Range("K132").Value = 19487.25
$ws.Range("M132").Value = -16957.25
$ws.Range("H136").Value = 6431.5557
$ws.Range("J136").Value = 6597.6
$ws.Range("L136").Value = 19792.8
$ws.Range("N136").Value = -24892.8

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 2669.182
$ws.Range("I132").Value = 2436.1
$ws.Range("J132").Value = 5000
$ws.Range("K132").Value = 7308.299999999999
$ws.Range("L132").Value = 15000
$ws.Range("M132").Value = -4778.299999999999
$ws.Range("N132").Value = -20060
